# Replace the hadron label "pi-" with "h-" throughout the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$used.Replace("pi-", "h-", 1, 1, $false, $false, $false, $false)
